$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 76: 2025-10-08, site 四方坪站充电量(kw) (shared string index 26 in source file)
$ws.Cells.Item(76,1).Value = 45938
$ws.Cells.Item(76,2).Value = "四方坪站充电量(kw)"
$ws.Cells.Item(76,3).Value = 660.81
$ws.Cells.Item(76,4).Value = 1503.4819999999995
$ws.Cells.Item(76,5).Value = 493.37
$ws.Cells.Item(76,6).Value = 575.68700000000001
$ws.Cells.Item(76,7).Value = 467.75500000000005
$ws.Cells.Item(76,8).Value = 961.72299999999996
$ws.Cells.Item(76,9).Value = 695.04700000000003
$ws.Cells.Item(76,10).Value = 233.27099999999999
$ws.Cells.Item(76,11).Value = 246.87
$ws.Cells.Item(76,12).Value = 263.73
$ws.Cells.Item(76,13).Value = 316.19799999999998
$ws.Cells.Item(76,14).Value = 239.143
$ws.Cells.Item(76,15).Value = 948.20999999999981
$ws.Cells.Item(76,16).Value = 1408.0909999999997
$ws.Cells.Item(76,17).Value = 1026.7720000000002
$ws.Cells.Item(76,18).Value = 411.98799999999994
$ws.Cells.Item(76,19).Value = 430.79299999999995
$ws.Cells.Item(76,20).Value = 144.74700000000001
$ws.Cells.Item(76,21).Value = 114.69499999999999
$ws.Cells.Item(76,22).Value = 43.379999999999995
$ws.Cells.Item(76,23).Value = 70.057999999999993
$ws.Cells.Item(76,24).Value = 179.95
$ws.Cells.Item(76,25).Value = 88.56
$ws.Cells.Item(76,26).Value = 35.08

# Row 77: 2025-10-08, site 高岭站充电量(kw) (shared string index 27 in source file)
$ws.Cells.Item(77,1).Value = 45938
$ws.Cells.Item(77,2).Value = "高岭站充电量(kw)"
$ws.Cells.Item(77,3).Value = 597.40899999999988
$ws.Cells.Item(77,4).Value = 658.54899999999998
$ws.Cells.Item(77,5).Value = 221.90199999999999
$ws.Cells.Item(77,6).Value = 27.780999999999999
$ws.Cells.Item(77,7).Value = 10.654999999999999
$ws.Cells.Item(77,8).Value = 245.54400000000001
$ws.Cells.Item(77,9).Value = 149.524
$ws.Cells.Item(77,10).Value = 171.054
$ws.Cells.Item(77,11).Value = 213.858
$ws.Cells.Item(77,12).Value = 390.09299999999996
$ws.Cells.Item(77,13).Value = 105.744
$ws.Cells.Item(77,14).Value = 356.48800000000006
$ws.Cells.Item(77,15).Value = 217.072
$ws.Cells.Item(77,16).Value = 447.46599999999995
$ws.Cells.Item(77,17).Value = 572.20499999999993
$ws.Cells.Item(77,18).Value = 414.30200000000002
$ws.Cells.Item(77,19).Value = 27.597999999999999
$ws.Cells.Item(77,20).Value = 204.79400000000001
$ws.Cells.Item(77,21).Value = 194.50200000000001
$ws.Cells.Item(77,22).Value = 82.32
$ws.Cells.Item(77,23).Value = 70.768999999999991
$ws.Cells.Item(77,24).Value = 153.13200000000001
$ws.Cells.Item(77,25).Value = 27.375
$ws.Cells.Item(77,26).Value = 19.997

# Move the selection to G82, matching the post-edit cursor position
$ws.Range("G82").Select()
